{"js": "// Remove the \"Ver no Jupiter ...\" line, the \"\u00a9 2020 ...\" footer line, and the\n// now-orphaned blank paragraph that separated them from the \"Requisitos\"\n// section above, matching the upstream Jekyll site rebuild that dropped the\n// page-footer boilerplate from this export.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\nlet footerStart = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targetTexts[0]) {\n    footerStart = i;\n    break;\n  }\n}\n\nif (footerStart === -1) {\n  throw new Error(\"Could not locate the 'Ver no Jupiter' paragraph to remove.\");\n}\n\n// The blank paragraph immediately preceding the footer block is the spacer\n// that was only there to separate the footer from \"Requisitos\"; it goes too.\nconst deleteIndices = [];\nif (footerStart - 1 >= 0 && items[footerStart - 1].text === \"\") {\n  deleteIndices.push(footerStart - 1);\n}\ndeleteIndices.push(footerStart);\nif (footerStart + 1 < items.length && items[footerStart + 1].text === targetTexts[1]) {\n  deleteIndices.push(footerStart + 1);\n}\n\n// Delete from the end backwards so earlier indices stay valid.\ndeleteIndices.sort((a, b) => b - a);\nfor (const idx of deleteIndices) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Ver no Jupiter ...\" line, the \"\u00a9 2020 ...\" footer line, and the\n# blank spacer paragraph immediately above them, matching the upstream Jekyll\n# site rebuild that dropped the page-footer boilerplate from this export.\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n$footerStart = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $txt = $paras.Item($i).Range.Text.TrimEnd()\n    if ($txt -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n        $footerStart = $i\n        break\n    }\n}\n\nif ($footerStart -eq -1) {\n    throw \"Could not locate the 'Ver no Jupiter' paragraph to remove.\"\n}\n\n$deleteIdx = @()\n\n# The blank paragraph right before the footer block was only there to\n# separate it from \"Requisitos\"; it goes too.\nif (($footerStart - 1) -ge 1 -and $paras.Item($footerStart - 1).Range.Text.TrimEnd() -eq \"\") {\n    $deleteIdx += ($footerStart - 1)\n}\n\n$deleteIdx += $footerStart\n\nif (($footerStart + 1) -le $count -and $paras.Item($footerStart + 1).Range.Text.TrimEnd() -eq \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\") {\n    $deleteIdx += ($footerStart + 1)\n}\n\n# Delete from the highest index down so earlier indices stay valid.\n$sortedIdx = $deleteIdx | Sort-Object -Descending\nforeach ($idx in $sortedIdx) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
